# "Add files via upload" — a new line was inserted into the Sheet1 command
# listing: "conda install -c menpo opencv3", placed right after the
# "source activate carnd-term1" row and before the Udacity sign-in URL row.
# This pushes every following entry down one shared-string slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 holds "source activate carnd-term1"; row 8 was blank and is where
# the new command line now lives (row 9 already held the next entry).
$ws.Range("A8").Value = "conda install -c menpo opencv3"

# Minor re-wrap height adjustments observed on the two multi-line / styled
# cells around the edit.
$ws.Rows.Item(5).RowHeight = 23.95
$ws.Rows.Item(17).RowHeight = 35.2

# Final cursor position left on the sheet after the edit.
$ws.Range("B16").Select()
